# Apply cryptos list update (prices and 1h volume % changes, plus a couple of
# coin re-rankings in rows 40/41 and row 51) as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.395.14'
$ws.Range("E2").Value = '  -0.12%  '
# Row 3
$ws.Range("D3").Value = '2.067.47'
$ws.Range("E3").Value = '  +0.05%  '
# Row 4
$ws.Range("E4").Value = '  +0.03%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.40%  '
# Row 7
$ws.Range("E7").Value = '  +0.06%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.17'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.96%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.42%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0770'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.95%  '
# Row 11
$ws.Range("E11").Value = '  +0.78%  '
# Row 12
$ws.Range("D12").Value = '2.372.90'
$ws.Range("E12").Value = '  +0.08%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.46%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.19%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.776'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.50%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.12%  '
# Row 17
$ws.Range("D17").Value = '2.069.07'
$ws.Range("E17").Value = '  -0.33%  '
# Row 18
$ws.Range("D18").Value = '37.351.03'
$ws.Range("E18").Value = '  -0.79%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.70%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.65'
$ws.Range("D20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = '0.0₃0815'
$ws.Range("E21").Value = '  +0.13%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.75%  '
# Row 23
$ws.Range("E23").Value = '  +0.00%  '
# Row 24
$ws.Range("E24").Value = '  +0.73%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.19%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.44%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.94%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.00%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.40%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.61%  '
# Row 31
$ws.Range("E31").Value = '  -1.17%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.36%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.92%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.67%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.74%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.78'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '
# Row 37
$ws.Range("E37").Value = '  -0.05%  '
# Row 38
$ws.Range("E38").Value = '  -2.47%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.96%  '
# Row 40
$ws.Range("B40").Value = 'HuobiToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.63%  '
# Row 41
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.01'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.92%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0965'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.65%  '
# Row 43
$ws.Range("D43").Value = '1.479.16'
$ws.Range("E43").Value = '  -0.45%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.44%  '
# Row 45
$ws.Range("E45").Value = '  +0.97%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.73%  '
# Row 47
$ws.Range("E47").Value = '  +0.01%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.81%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.40%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.39%  '
# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.262.78'
$ws.Range("E51").Value = '  +0.16%  '
